$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Flip the formula in column F from D/E to E/D for rows 2 through 30.
# F2 holds its own (non-shared) formula; F3:F30 form a shared formula
# group. Set them separately so the original grouping is preserved.
$ws.Range("F2").Formula = "=E2/D2"
$ws.Range("F3:F30").Formula = "=E3/D3"

# Update the selection to match the new active cell / selected range.
$ws.Range("F2:F30").Select()
